# Global_Parameters.xlsx - add a "Mip Gap" parameter block to the
# "Solver Options" section of the "Global Parameters" sheet (between the
# existing "pEnableRMIP" row and the "Scaling" section), and extend the
# pEnableRMIP Yes/No dropdown down into the newly blank row below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global Parameters")

# ---------------------------------------------------------------------
# 1) Make room: insert 3 new rows at row 10 (pushes "Scaling" section and
#    everything below it down by 3 rows, e.g. old row 10 -> 13, old row
#    19 -> 22).
# ---------------------------------------------------------------------
$ws.Rows("10:12").Insert()

# ---------------------------------------------------------------------
# 2) Copy cell formatting from analogous existing rows so the new rows
#    pick up the same styles (fonts/fills/borders/number formats) used
#    elsewhere on the sheet, rather than plain defaults.
# ---------------------------------------------------------------------

# Row 10 ("Mip Gap" sub-header) should look like row 7 ("Solve as rMIP"):
# bold label in B, plain value box in C, spacer formatting in E/F/H.
$ws.Range("B7:C7").Copy()
$ws.Range("B10:C10").PasteSpecial(-4122)
$ws.Range("E7:F7").Copy()
$ws.Range("E10:F10").PasteSpecial(-4122)
$ws.Range("H7").Copy()
$ws.Range("H10").PasteSpecial(-4122)

# Row 11 (pMIPGap value row) should look like the pMovWindow row, which
# is now row 22 (old row 19) after the insert above - same layout of a
# numeric parameter row with a unit label in G.
$ws.Range("B22:C22").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)
$ws.Range("E22:H22").Copy()
$ws.Range("E11:H11").PasteSpecial(-4122)

# Row 9 is blank below pEnableRMIP; extend its B/C formatting (currently
# only E/F/H are formatted) to match the "label" style used in column B,
# since the Yes/No dropdown area now spans C8:C9.
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Row heights: give the two new content rows the same explicit row
#    height as the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Rows("10").RowHeight = 18.75
$ws.Rows("11").RowHeight = 18.75

# ---------------------------------------------------------------------
# 4) Fill in the new "Mip Gap" content.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Mip Gap"
$ws.Range("C10").Value = "[%]"

$ws.Range("B11").Value = "pMIPGap"
$ws.Range("C11").Value = 0.05
$ws.Range("E11").Value = "Relative MIP gap"
$ws.Range("F11").Value = "The MIP solver will terminate (with an optimal result) when the gap between the lower and upper objective bound is less than pMIPGap"
$ws.Range("G11").Value = "Factor"
$ws.Range("H11").Value = 0.05

# ---------------------------------------------------------------------
# 5) Conditional formatting: pEnableRMIP's "No"/"Yes" highlight now
#    covers C8:C9 as well, and the new pMIPGap cell gets the same kind
#    of (otherwise unused, since it is numeric) No/Yes rule so it keeps
#    parity with the other parameter cells on the sheet.
# ---------------------------------------------------------------------
$ws.Range("C9").FormatConditions.Delete()
$fc = $ws.Range("C9").FormatConditions.Add(1, 3, '="No"')
$fc.Font.Bold = $true
$fc.Font.Color = 3473849
$fc = $ws.Range("C9").FormatConditions.Add(1, 3, '="Yes"')
$fc.Font.Bold = $true
$fc.Font.Color = 4824142

$ws.Range("C11").FormatConditions.Delete()
$fc = $ws.Range("C11").FormatConditions.Add(1, 3, '="No"')
$fc.Font.Bold = $true
$fc.Font.Color = 3473849
$fc = $ws.Range("C11").FormatConditions.Add(1, 3, '="Yes"')
$fc.Font.Bold = $true
$fc.Font.Color = 4824142

# ---------------------------------------------------------------------
# 6) Data validation:
#      - pEnableRMIP's Yes/No list now covers C8:C9
#      - the free-input validation group gains C11 (pMIPGap)
# ---------------------------------------------------------------------
$ws.Range("C8:C9").Validation.Delete()
$ws.Range("C8:C9").Validation.Add(3, 1, 1, "No, Yes")
$ws.Range("C8:C9").Validation.InputMessage = ""
$ws.Range("C8:C9").Validation.ErrorMessage = ""

$ws.Range("C5,C15,C18,C11").Validation.Delete()
$ws.Range("C5").Validation.Add(0, 1, 1)
$ws.Range("C15").Validation.Add(0, 1, 1)
$ws.Range("C18").Validation.Add(0, 1, 1)
$ws.Range("C11").Validation.Add(0, 1, 1)

Write-Host "Mip Gap block inserted."
